$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "Denomination"
$ws.Range("I2").Value = "Y"
$ws.Range("I3").Value = "Y"
$ws.Range("I4").Value = "Y"
$ws.Range("I8").Value = "Y"
$ws.Range("I9").Value = "y"
$ws.Range("I10").Value = "Y"

$ws.Range("I10").Select()
